# Applies updated crypto price/volume data as described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellAddress, $NewValue)
    $range = $ws.Range($CellAddress)
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    $range.Style = "Normal"
}

Set-TextValue "D2" '67.297.91'

Set-TextValue "D3" '3.113.77'
Set-TextValue "E3" '  -0.87%  '

Set-TextValue "E4" '  -0.02%  '

Set-TextValue "D5" '579.94'
Set-TextValue "E5" '  -0.21%  '

Set-TextValue "D6" '173.17'
Set-TextValue "E6" '  -0.82%  '

Set-TextValue "D7" '0.999'
Set-TextValue "E7" '  -0.07%  '

Set-TextValue "E8" '  -1.02%  '

Set-TextValue "D9" '6.50'
Set-TextValue "E9" '  +0.14%  '

Set-TextValue "E10" '  -1.83%  '

Set-TextValue "E11" '  -1.69%  '

Set-TextValue "E12" '  -1.51%  '

Set-TextValue "E13" '  -1.78%  '

Set-TextValue "E14" '  -1.18%  '

Set-TextValue "D15" '3.629.02'
Set-TextValue "E15" '  -0.83%  '

Set-TextValue "D16" '67.204.19'
Set-TextValue "E16" '  -0.06%  '

Set-TextValue "D17" '7.10'
Set-TextValue "E17" '  -1.48%  '

Set-TextValue "D18" '3.110.60'
Set-TextValue "E18" '  -0.88%  '

Set-TextValue "D19" '16.61'
Set-TextValue "E19" '  +2.29%  '

Set-TextValue "D20" '491.11'
Set-TextValue "E20" '  +0.53%  '

Set-TextValue "D21" '0.700'
Set-TextValue "E21" '  -2.76%  '

Set-TextValue "D22" '7.82'
Set-TextValue "E22" '  +1.59%  '

Set-TextValue "D23" '83.85'
Set-TextValue "E23" '  -0.78%  '

Set-TextValue "D24" '13.10'
Set-TextValue "E24" '  -2.24%  '

Set-TextValue "E25" '  -2.21%  '

Set-TextValue "D26" '10.57'
Set-TextValue "E26" '  +4.47%  '

Set-TextValue "E27" '  -0.06%  '

Set-TextValue "E29" '  -2.86%  '

Set-TextValue "E30" '  -1.36%  '

Set-TextValue "D31" '28.27'
Set-TextValue "E31" '  -2.77%  '

Set-TextValue "E32" '  -1.53%  '

Set-TextValue "D33" '0.0₃0945'
Set-TextValue "E33" '  -6.29%  '

Set-TextValue "E34" '  -0.07%  '

Set-TextValue "D35" '5.79'
Set-TextValue "E35" '  -3.15%  '

Set-TextValue "E36" '  -2.34%  '

Set-TextValue "D37" '46.73'
Set-TextValue "E37" '  -1.70%  '

Set-TextValue "E38" '  -4.45%  '

Set-TextValue "E39" '  -0.05%  '

Set-TextValue "E40" '  -2.39%  '

Set-TextValue "D41" '8.46'
Set-TextValue "E41" '  -2.67%  '

Set-TextValue "D42" '385.22'
Set-TextValue "E42" '  -1.12%  '

Set-TextValue "D43" '2.800.60'
Set-TextValue "E43" '  -2.27%  '

Set-TextValue "E44" '  -9.17%  '

Set-TextValue "D45" '0.0350'
Set-TextValue "E45" '  -2.82%  '

Set-TextValue "D46" '135.31'
Set-TextValue "E46" '  -0.94%  '

Set-TextValue "E47" '  -0.01%  '

Set-TextValue "D48" '24.99'
Set-TextValue "E48" '  -1.06%  '

Set-TextValue "E49" '  -1.88%  '

Set-TextValue "E50" '  -1.99%  '

Set-TextValue "D51" '6.70'
Set-TextValue "E51" '  -2.01%  '

Write-Output "Updated cryptos list applied."